$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.331000000000007
$ws.Range("A9").Value = -22.09589999999999
$ws.Range("B9").Value = 6.482800000000005
$ws.Range("C9").Value = -11.9614
$ws.Range("B11").Value = 6.1548
$ws.Range("A18").Value = -22.07240000000001
$ws.Range("A20").Value = -21.57149999999997
$ws.Range("B23").Value = 8.806599999999994
$ws.Range("B24").Value = 5.304600000000002
$ws.Range("B26").Value = 5.765800000000002
$ws.Range("A27").Value = -21.83259999999999
$ws.Range("C27").Value = -12.09529999999999
$ws.Range("C29").Value = -11.7221
$ws.Range("C32").Value = -12.80590000000001
$ws.Range("B34").Value = 9.930700000000007
$ws.Range("A35").Value = -21.9699
$ws.Range("B35").Value = 4.9359
$ws.Range("C37").Value = -12.91499999999999
$ws.Range("C38").Value = -11.59180000000001
$ws.Range("C41").Value = -12.96460000000001
$ws.Range("C45").Value = -14.07449999999999
$ws.Range("B48").Value = 5.652100000000005
$ws.Range("B49").Value = 5.732199999999999
$ws.Range("C51").Value = -11.8571
$ws.Range("B52").Value = 5.605999999999998
$ws.Range("C57").Value = -14.15089999999999
$ws.Range("C64").Value = -10.19049999999999
$ws.Range("B66").Value = 5.406399999999996
$ws.Range("B67").Value = 5.403999999999997
$ws.Range("A69").Value = -21.63709999999999
$ws.Range("A76").Value = -19.6372
$ws.Range("A78").Value = -21.54509999999999
$ws.Range("B78").Value = 5.594099999999997
$ws.Range("B80").Value = 9.712200000000003
$ws.Range("A82").Value = -21.83869999999999
$ws.Range("C82").Value = -11.349
$ws.Range("A83").Value = -21.53539999999998
$ws.Range("A93").Value = -21.43630000000001
$ws.Range("C93").Value = -10.50639999999999
$ws.Range("B99").Value = 5.628400000000002
$ws.Range("C102").Value = -11.844
$ws.Range("B104").Value = 10.0205
$ws.Range("C105").Value = -12.59510000000001
